$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on B:D columns so numeric-looking strings (e.g. "1.001", "108.00")
# are preserved verbatim as text instead of being coerced to Double and losing formatting.
$ws.Range("B2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.487.78'
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '1.731.70'
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '247.13'
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '0.4883'
$ws.Range("E7").Value = '  +1.19%  '
$ws.Range("D8").Value = '0.2669'
$ws.Range("E8").Value = '  -0.87%  '
$ws.Range("D9").Value = '0.06219'
$ws.Range("E9").Value = '  -0.70%  '
$ws.Range("D10").Value = '1.733.28'
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("D11").Value = '0.07068'
$ws.Range("E11").Value = '  -1.06%  '
$ws.Range("E12").Value = '  -1.65%  '
$ws.Range("D13").Value = '4.660'
$ws.Range("E13").Value = '  +3.05%  '
$ws.Range("D14").Value = '0.6095'
$ws.Range("E14").Value = '  -2.48%  '
$ws.Range("D15").Value = '77.45'
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("D17").Value = '26.483.23'
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").Value = '0.000007147'
$ws.Range("E19").Value = '  +3.44%  '
$ws.Range("E20").Value = '  -2.50%  '
$ws.Range("D21").Value = '1.958.41'
$ws.Range("E21").Value = '  -0.48%  '
$ws.Range("D22").Value = '4.525'
$ws.Range("E22").Value = '  -2.10%  '
$ws.Range("E23").Value = '  -1.18%  '
$ws.Range("D24").Value = '5.257'
$ws.Range("D25").Value = '139.52'
$ws.Range("E25").Value = '  +2.39%  '
$ws.Range("D26").Value = '15.42'
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("D27").Value = '1.775'
$ws.Range("E27").Value = '  -2.22%  '
$ws.Range("D28").Value = '108.00'
$ws.Range("E28").Value = '  +1.11%  '
$ws.Range("E29").Value = '  -2.00%  '
$ws.Range("D30").Value = '3.971'
$ws.Range("E30").Value = '  -1.12%  '
$ws.Range("E31").Value = '  +1.86%  '
$ws.Range("D32").Value = '3.692'
$ws.Range("E32").Value = '  -1.44%  '
$ws.Range("D33").Value = '0.04581'
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '2.615'
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.004'
$ws.Range("E35").Value = '  +0.42%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.6365'
$ws.Range("E36").Value = '  -0.84%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '0.8990'
$ws.Range("E37").Value = '  -3.50%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '2.019'
$ws.Range("E38").Value = '  +1.68%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.399'
$ws.Range("E39").Value = '  -1.59%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '1.002'
$ws.Range("E40").Value = '  +0.23%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.01508'
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '101.37'
$ws.Range("E42").Value = '  -11.05%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.447'
$ws.Range("E43").Value = '  -6.05%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.3884'
$ws.Range("E44").Value = '  -1.03%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '6.947'
$ws.Range("E45").Value = '  +2.82%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.1184'
$ws.Range("E46").Value = '  -2.80%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.05384'
$ws.Range("E47").Value = '  +0.87%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = '30.58'
$ws.Range("E48").Value = '  -0.74%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '7.814'
$ws.Range("E49").Value = '  -1.50%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.248'
$ws.Range("E50").Value = '  -1.35%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '0.3412'
$ws.Range("E51").Value = '  -1.26%  '
